# Corrected excel sheets for application fix issues
#
# Applies the small rounding corrections that ripple through the
# "Summary", "Repayment schedule" and "Transactions" sheets after the
# underlying loan-schedule numbers were recalculated, and leaves the
# workbook with the "Transactions" tab active/selected (matching the
# last sheet the author was looking at when they saved).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = 836.76
$wsSummary.Range("E2").Value = 9163.24

$wsSummary.Range("A3").Value = 561.21
$wsSummary.Range("E3").Value = 510.25

$wsSummary.Range("A5").Value = 0.89
$wsSummary.Range("B5").Value = 0.51

# ---------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

$wsSchedule.Range("G3").Value = 4163.24

$wsSchedule.Range("G5").Value = 8313.84

$wsSchedule.Range("F6").Value = 794.33
$wsSchedule.Range("G6").Value = 7519.51
$wsSchedule.Range("H6").Value = 93.39

$wsSchedule.Range("F7").Value = 813.55
$wsSchedule.Range("G7").Value = 6705.96
$wsSchedule.Range("H7").Value = 74.17

$wsSchedule.Range("F8").Value = 819.37
$wsSchedule.Range("G8").Value = 5886.59
$wsSchedule.Range("H8").Value = 68.35

$wsSchedule.Range("F9").Value = 829.66
$wsSchedule.Range("G9").Value = 5056.93
$wsSchedule.Range("H9").Value = 58.06

$wsSchedule.Range("F10").Value = 836.18
$wsSchedule.Range("G10").Value = 4220.75
$wsSchedule.Range("H10").Value = 51.54

$wsSchedule.Range("F11").Value = 844.7
$wsSchedule.Range("G11").Value = 3376.05
$wsSchedule.Range("H11").Value = 43.02

$wsSchedule.Range("F12").Value = 854.42
$wsSchedule.Range("G12").Value = 2521.63
$wsSchedule.Range("H12").Value = 33.3

$wsSchedule.Range("F13").Value = 862.02
$wsSchedule.Range("G13").Value = 1659.61
$wsSchedule.Range("H13").Value = 25.7

$wsSchedule.Range("F14").Value = 871.35
$wsSchedule.Range("G14").Value = 788.26
$wsSchedule.Range("H14").Value = 16.37

$wsSchedule.Range("F15").Value = 788.26

$wsSchedule.Range("K15").Value = 796.29
$wsSchedule.Range("P15").Value = 796.29

# ---------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")

$wsTransactions.Range("A2").Value = 3241
$wsTransactions.Range("J2").NumberFormat = "#,##0.00"
$wsTransactions.Range("J2").Value = 9163.24

$wsTransactions.Range("A3").Value = 3239
$wsTransactions.Range("J3").NumberFormat = "#,##0.00"
$wsTransactions.Range("J3").Value = 4163.24

$wsTransactions.Range("A4").Value = 3233

# ---------------------------------------------------------------------
# Restore the per-sheet selections exactly as they were left, switching
# through each sheet so every sheetView keeps its own last selection,
# and finish on "Transactions" so it becomes the active tab.
# ---------------------------------------------------------------------
$wsNewLoanInput = $wb.Worksheets.Item("NewLoanInput")
$wsNewLoanInput.Activate()
$wsNewLoanInput.Range("B2").Select() | Out-Null

$wsSummary.Activate()
$wsSummary.Range("D5").Select() | Out-Null

$wsSchedule.Activate()
$wsSchedule.Range("F13").Select() | Out-Null

$wsTransactions.Activate()
$wsTransactions.Range("H4").Select() | Out-Null
